# Renderer017-TemplateFormula: add a second worksheet demonstrating an
# incorrect variable path inside a TemplateFormulaCell, alongside the
# existing correct example.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the existing sheet.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Correct Formula"

# Move the selection the way the "after" workbook has it.
$ws1.Range("C6").Select()

# ---------------------------------------------------------------------
# 2) Add the new sheet right after the first one and populate it with
#    the same table, but pointed at an incorrect/undefined range.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Incorrect Formula"

# Column widths / row heights to match the template layout.
$ws2.Range("A1").ColumnWidth = 18.830729166666668
$ws2.Range("B1").ColumnWidth = 17.830729166666668
$ws2.Range("C1").ColumnWidth = 50

# Header row.
$ws2.Range("A1").Value = "Item"
$ws2.Range("B1").Value = "Weight"
$ws2.Range("C1").Value = "Price"
$headerRange = $ws2.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12

# Data rows.
$ws2.Range("A2").Value = "Item A"
$ws2.Range("B2").Value = 13.5
$ws2.Range("C2").Value = 12.1

$ws2.Range("A3").Value = "Item B"
$ws2.Range("B3").Value = 13.5
$ws2.Range("C3").Value = 12

$ws2.Range("A4").Value = "Item C"
$ws2.Range("B4").Value = 3.5
$ws2.Range("C4").Value = 12.7

$ws2.Range("A5").Value = "Item D"
$ws2.Range("B5").Value = 13
$ws2.Range("C5").Value = 2.1

$ws2.Range("A6").Value = "Item E"
$ws2.Range("B6").Value = 1.5
$ws2.Range("C6").Value = 32

# Footer / custom-formula row, styled with the "incorrect" (orange) fill.
$ws2.Range("A8:C8").Interior.Pattern = 1
$ws2.Range("A8:C8").Interior.ThemeColor = 6
$ws2.Range("A8:C8").Font.Size = 12
$ws2.Range("A8:C8").Font.Color = $ws2.Range("A8").Font.Color

$ws2.Range("B8").Value = "Incorrect custom formula:"
$ws2.Range("B8").Font.Bold = $false
$ws2.Range("B8").HorizontalAlignment = -4152
$ws2.Range("B8").VerticalAlignment = -4108

$ws2.Range("C8").Formula = "=MAX( Bundefined:Bundefined)"
$ws2.Range("C8").Font.Bold = $true
$ws2.Range("C8").HorizontalAlignment = -4108
$ws2.Range("C8").VerticalAlignment = -4108
$ws2.Range("C8").WrapText = $true

$ws2.Range("A8").RowHeight = 51

$ws2.Range("C14").Select()
